$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: LinearRegression values update
$ws.Range("B2").Value = 309673752764775.1
$ws.Range("C2").Value = 309673752764775.1
$ws.Range("D2").Value = 309673752764775.1

# Row 3: RandomForestRegressor values update
$ws.Range("B3").Value = 5353635546107.014
$ws.Range("C3").Value = 5082660079774.208
$ws.Range("D3").Value = 6785061732764.529

# Row 4: model name changed from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 4016399359364.119
$ws.Range("C4").Value = 4228732044479.397
$ws.Range("D4").Value = 4150422191803.854

# Row 5: model name changed from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 28918929791744.3
$ws.Range("C5").Value = 87324136894645.14
$ws.Range("D5").Value = 83651572154793.67
